$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New air hubs: Malmo (row 17) and Hong Kong (row 18) - "seed for gc air"
# ---------------------------------------------------------------------------

# -- Row 17: Malmo (Sweden) --------------------------------------------------
$ws.Range("A17").Value = "active"
$ws.Range("B17").Value = "air"
$ws.Range("C17").Value = "Malmo"
$ws.Range("D17").Value = "MMX"
$ws.Range("F17").Value = 55.535558
$ws.Range("G17").Value = 13.363027
$ws.Range("H17").Value = "Sweden"
$ws.Range("I17").Value = "230 32 Malmö-Sturup, Sweden"
$ws.Range("J17").Value = "zipcode"

# -- Row 18: Hong Kong (China) -----------------------------------------------
$ws.Range("A18").Value = "active"
$ws.Range("B18").Value = "air"
$ws.Range("C18").Value = "Hong Kong"
$ws.Range("D18").Value = "HKG"
$ws.Range("F18").Value = 22.316265
$ws.Range("G18").Value = 113.939724
$ws.Range("H18").Value = "China"
$ws.Range("I18").Value = "1 Sky Plaza Rd, Chek Lap Kok, Hong Kong"

# ---------------------------------------------------------------------------
# Formatting - reuse the existing "plain" cell style (as seen on the rest of
# the table, e.g. C16) for most of the new cells.
# ---------------------------------------------------------------------------
$ws.Range("C16").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("H18").PasteSpecial(-4122)

# G18 (LONGITUDE) reuses the fill+font combo already used throughout column G
$ws.Range("G2").Copy()
$ws.Range("G18").PasteSpecial(-4122)

# Build each *new* look once on a scratch cell, then stamp it onto every
# target cell with a single PasteSpecial so no unused transient styles are
# left behind in the stylesheet.

# -- Roboto, left aligned (F17, G17, F18) ------------------------------------
$ws.Range("Z100").Value = "x"
$ws.Range("G2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("Z100").Font.Name = "Roboto"
$ws.Range("Z100").HorizontalAlignment = -4131
$ws.Range("Z100").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# -- Dark grey Arial, left aligned (I17) -------------------------------------
$ws.Range("Z100").Value = "x"
$ws.Range("G2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("Z100").Font.Color = 2236962
$ws.Range("Z100").HorizontalAlignment = -4131
$ws.Range("Z100").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# I18 reuses the plain style too
$ws.Range("C16").Copy()
$ws.Range("I18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column widths: split the former shared 9:10 width so FULL_ADDRESS (col I)
# is wider than TRUCKING_METHOD (col J).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 44.29
$ws.Columns.Item(10).ColumnWidth = 28

Write-Output "done"
